# Season record columns (Wins/Losses/Ties) were missing from the team
# statistics export - add them now, pulled from the season record data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, formatted like the rest of the header row (copy the
# style from the last existing header cell, AC1, which carries the bold /
# bordered / centered header style).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row shares the same team season record: 63-99-0.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 63
    $ws.Cells.Item($r, 31).Value = 99
    $ws.Cells.Item($r, 32).Value = 0
}
